$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("G8").Value = 2.27
$ws.Range("I8").Value = 2.85
$ws.Range("J8").Value = 2.85
$ws.Range("K8").Value = 2.15
$ws.Range("N8").Value = 7.8
$ws.Range("P8").Value = 3.6
$ws.Range("Q8").Value = 1.75
$ws.Range("R8").Value = 2
$ws.Range("T8").Value = 2.87
$ws.Range("Z8").Value = 24
$ws.Range("AB8").Value = 23
$ws.Range("AC8").Value = 7.8
$ws.Range("AD8").Value = 6.5
$ws.Range("AG8").Value = 10.5
$ws.Range("AN8").Value = 4.35
$ws.Range("AP8").Value = 18
$ws.Range("AQ8").Value = 45
$ws.Range("AR8").Value = 70
$ws.Range("AT8").Value = 2.87
